$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated prediction values (X_Pred, Y_Pred, Z_Pred) for rows 2-102
# Each entry: row, E (X_Pred), F (Y_Pred), G (Z_Pred)
$data = @(
  @(2, -9.915353532597596, -3.058207550032201, -6.72473314546327),
  @(3, -10.20825453177269, -2.99638569601357, -6.691373957822552),
  @(4, -11.03263756480766, -2.796937554522716, -6.568384864927377),
  @(5, -11.76034703365932, -2.77292627111099, -6.81438923532341),
  @(6, -12.62039349963514, -2.568987469745423, -6.319513280209006),
  @(7, -13.30188713945551, -2.505463616357368, -6.405307140730837),
  @(8, -13.87788300297725, -2.416396680125106, -5.887310178797946),
  @(9, -14.88473346071735, -2.238197346146372, -5.673813996357918),
  @(10, -15.63178026086624, -2.205689158190366, -5.296768766819049),
  @(11, -16.34161873633894, -1.996107574300015, -5.061644100084491),
  @(12, -16.99183486467043, -2.129243201897528, -4.541827308056599),
  @(13, -17.96252438196168, -1.912225189993076, -3.859875437636771),
  @(14, -18.55154708681094, -1.659373545210821, -3.599718287868833),
  @(15, -19.48591855601924, -1.52481085660356, -3.008404430022275),
  @(16, -20.40959361380586, -1.350369013540402, -2.570768022931898),
  @(17, -21.20649590317604, -1.183429060005553, -2.085960048702781),
  @(18, -22.16477391737334, -1.07948926774509, -1.704005205598265),
  @(19, -23.01149241905629, -0.9919410386424676, -1.363919546981589),
  @(20, -23.71144620588349, -0.7859991149421166, -1.025313318586027),
  @(21, -24.32053941098911, -0.5221630280756528, -0.656267486083453),
  @(22, -24.998432667528, -0.3868278936007305, -0.4424309037695411),
  @(23, -25.40976663820935, -0.3256606547241863, -0.2496729290308082),
  @(24, -25.85830893356689, -0.3620048874128179, -0.2098068668777523),
  @(25, -26.10850284087228, -0.05814563075919755, -0.1960337642882565),
  @(26, -26.12399103513404, -0.09735707777016712, -0.3602243342263431),
  @(27, -26.14821179539123, -0.291502836610123, -0.3311856065233947),
  @(28, -26.18837898050965, -0.07868745391786575, -0.6319027104950008),
  @(29, -26.29488486412709, -0.1340417103326836, -0.6929652109488111),
  @(30, -26.02568093309552, -0.2337265041695858, -0.8309318982949204),
  @(31, -25.60714619585122, -0.2423674240451251, -1.354584735055438),
  @(32, -25.68795188899035, -0.1546228103998769, -1.528882562787337),
  @(33, -25.25983358606591, -0.2005898856771772, -2.034179990963755),
  @(34, -24.62048406909306, -0.2411629321836863, -2.060835919549509),
  @(35, -24.26023626410012, -0.439472043327311, -2.208307618758711),
  @(36, -23.46270554419354, -0.5324535781092495, -2.375705802893021),
  @(37, -22.66812059242634, -0.4844833804971652, -2.543051617815963),
  @(38, -22.29715019140603, -0.5424822820860117, -2.817296085441599),
  @(39, -21.63804438944502, -0.4872851333052946, -3.010708675322419),
  @(40, -21.23878152198374, -0.6620673762423376, -3.201096943246799),
  @(41, -20.60190645024797, -0.4630905576537848, -3.117293112756911),
  @(42, -20.10813024857228, -0.2156984031565291, -3.234050269499424),
  @(43, -19.18620955936645, -0.1882045671889043, -3.217724167855792),
  @(44, -18.85064074523017, -0.1368041862322878, -3.45246915780794),
  @(45, -18.27419974341181, -0.1646776989823226, -3.524057869746498),
  @(46, -17.80834942369752, -0.08787825051275745, -3.33270077141183),
  @(47, -17.11977284803978, -0.1441751527321795, -3.716842029090914),
  @(48, -16.46060158456455, -0.09352103303754139, -3.568375314865741),
  @(49, -15.80868345686364, -0.005174173461573876, -3.874054401614361),
  @(50, -15.08556938630943, 0.06471253910755988, -3.836361661733032),
  @(51, -14.26709098185608, 0.09417022050144352, -4.015477456910686),
  @(52, -13.68778276571539, 0.1631404718716565, -3.911799510707057),
  @(53, -13.01041320129016, 0.2374785674069771, -3.970767242706192),
  @(54, -12.62569588228604, 0.09588531217370964, -4.176879366343485),
  @(55, -12.02262513648762, -0.05322292489070855, -4.58453439992631),
  @(56, -11.59199311141756, -0.07312322521013218, -4.709225492190908),
  @(57, -10.81588139896004, -0.114665102126929, -5.137068856755675),
  @(58, -10.55924607865653, -0.275634965565951, -5.142921116125926),
  @(59, -9.936013186481841, -0.2276778602567084, -5.455067800478359),
  @(60, -9.627021747114263, -0.3297323609079631, -5.327627324616998),
  @(61, -9.05782078876738, -0.3561133511340411, -5.698702464060039),
  @(62, -8.583316456874702, -0.3059698312502303, -5.728186330059606),
  @(63, -8.131579639623785, -0.3165876888548702, -5.709608352327197),
  @(64, -7.791048842710491, -0.4660756027016984, -5.881798319301579),
  @(65, -7.434728728570074, -0.3882418623076369, -6.118415508560094),
  @(66, -7.178761115711488, -0.5119379395562648, -6.064475220852183),
  @(67, -7.014203961293833, -0.6072891811525557, -6.145267821688474),
  @(68, -6.732562342562622, -0.5593844450546801, -6.102704745150023),
  @(69, -6.626501597241799, -0.8585697495938042, -6.115116248243979),
  @(70, -6.443196265154793, -0.7558999107089883, -5.939901959313159),
  @(71, -6.370926753468464, -0.8272791458020788, -5.840806319104134),
  @(72, -6.392895637636881, -0.9635438337787638, -5.637823255846014),
  @(73, -6.451352769825188, -1.019552705335668, -5.468709980039438),
  @(74, -6.668724273906365, -0.898632196289486, -5.312793745497322),
  @(75, -6.83987994895625, -1.251757788536521, -5.142659270069092),
  @(76, -7.304957822802885, -1.190656011174186, -4.90310940497403),
  @(77, -7.722104775945962, -1.262428015352528, -4.84356561164986),
  @(78, -8.323905568368742, -1.304703061228462, -4.493830925838831),
  @(79, -8.733707739617611, -1.33192195883641, -4.358024468461607),
  @(80, -9.178440174848202, -1.438205273305542, -4.004545384037844),
  @(81, -10.05514005003871, -1.429957122515255, -4.019994301391081),
  @(82, -10.79234143845062, -1.516418690482014, -3.810884040403033),
  @(83, -11.57065265778555, -1.652094224830821, -3.488708652073838),
  @(84, -12.76685709152553, -1.716613093234847, -3.339809891854887),
  @(85, -13.63386866261316, -1.773237303025313, -3.005223000431736),
  @(86, -14.73825677652412, -1.957393634797031, -2.753916247384804),
  @(87, -15.78361169692174, -1.988697330891599, -2.321137084648707),
  @(88, -17.17218133631521, -1.91778941870081, -2.055651367624185),
  @(89, -18.78737873789895, -2.142950842972815, -1.874205142540704),
  @(90, -20.35825869436345, -2.347269321120792, -1.589067878950752),
  @(91, -22.03102295154227, -2.491965452127548, -1.724285182700099),
  @(92, -23.66495543849252, -2.609887823822975, -1.522297134457949),
  @(93, -25.76537976799632, -2.70860378724959, -1.493585714326044),
  @(94, -27.97158990446126, -2.551561614663086, -1.578829698128522),
  @(95, -30.22342671632962, -2.985100130964002, -1.940412918011313),
  @(96, -32.45438130516543, -2.982599501121233, -2.203908605003892),
  @(97, -34.70031348845217, -3.249564648366869, -2.372655296330899),
  @(98, -37.19894020888706, -3.503961184884449, -2.829642127021347),
  @(99, -39.636020013663, -3.447323882791141, -3.309056072479672),
  @(100, -42.20728283026379, -3.624344909514119, -3.590409660548366),
  @(101, -44.5104676538774, -3.830090448671845, -4.091164059638705),
  @(102, -46.83279342424562, -3.955239771535904, -4.401713483044446)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 5).Value = $row[1]
  $ws.Cells.Item($r, 6).Value = $row[2]
  $ws.Cells.Item($r, 7).Value = $row[3]
}
